$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44532
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(2, 11).Value = 2000
$ws.Cells.Item(2, 12).Value = 2200
$ws.Cells.Item(2, 13).Value = 2100
$ws.Cells.Item(2, 16).Value = 2100

$ws.Cells.Item(3, 4).Value = 44761
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 200
$ws.Cells.Item(3, 11).Value = 700
$ws.Cells.Item(3, 12).Value = 800
$ws.Cells.Item(3, 13).Value = 750
$ws.Cells.Item(3, 16).Value = 750

$ws.Cells.Item(4, 4).Value = 44761
$ws.Cells.Item(4, 9).Value = "Segunda"
$ws.Cells.Item(4, 10).Value = 150
$ws.Cells.Item(4, 11).Value = 600
$ws.Cells.Item(4, 12).Value = 600
$ws.Cells.Item(4, 13).Value = 600
$ws.Cells.Item(4, 16).Value = 600

$ws.Cells.Item(5, 4).Value = 45033
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 300
$ws.Cells.Item(5, 11).Value = 1500
$ws.Cells.Item(5, 12).Value = 1500
$ws.Cells.Item(5, 13).Value = 1500
$ws.Cells.Item(5, 16).Value = 1500

$ws.Cells.Item(6, 4).Value = 45020
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 1400
$ws.Cells.Item(6, 12).Value = 1500
$ws.Cells.Item(6, 13).Value = 1475
$ws.Cells.Item(6, 16).Value = 1475

$ws.Cells.Item(7, 4).Value = 45002
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 1200
$ws.Cells.Item(7, 12).Value = 1200
$ws.Cells.Item(7, 13).Value = 1200
$ws.Cells.Item(7, 16).Value = 1200

$ws.Cells.Item(8, 4).Value = 44868
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 200
$ws.Cells.Item(8, 11).Value = 700
$ws.Cells.Item(8, 12).Value = 800
$ws.Cells.Item(8, 13).Value = 750
$ws.Cells.Item(8, 16).Value = 750

$ws.Cells.Item(9, 4).Value = 45001
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 150
$ws.Cells.Item(9, 11).Value = 1300
$ws.Cells.Item(9, 12).Value = 1300
$ws.Cells.Item(9, 13).Value = 1300
$ws.Cells.Item(9, 16).Value = 1300

$ws.Cells.Item(10, 4).Value = 44999
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 200
$ws.Cells.Item(10, 11).Value = 1500
$ws.Cells.Item(10, 12).Value = 1500
$ws.Cells.Item(10, 13).Value = 1500
$ws.Cells.Item(10, 16).Value = 1500

$ws.Cells.Item(11, 4).Value = 45037
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 1500
$ws.Cells.Item(11, 12).Value = 1500
$ws.Cells.Item(11, 13).Value = 1500
$ws.Cells.Item(11, 16).Value = 1500

$ws.Cells.Item(12, 4).Value = 44797
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 240
$ws.Cells.Item(12, 11).Value = 750
$ws.Cells.Item(12, 12).Value = 850
$ws.Cells.Item(12, 13).Value = 800
$ws.Cells.Item(12, 16).Value = 800

$ws.Cells.Item(13, 4).Value = 44797
$ws.Cells.Item(13, 9).Value = "Segunda"
$ws.Cells.Item(13, 10).Value = 200
$ws.Cells.Item(13, 11).Value = 650
$ws.Cells.Item(13, 12).Value = 650
$ws.Cells.Item(13, 13).Value = 650
$ws.Cells.Item(13, 16).Value = 650

$ws.Cells.Item(14, 4).Value = 45014
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 60
$ws.Cells.Item(14, 11).Value = 1500
$ws.Cells.Item(14, 12).Value = 1500
$ws.Cells.Item(14, 13).Value = 1500
$ws.Cells.Item(14, 16).Value = 1500

$ws.Cells.Item(15, 4).Value = 44764
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 700
$ws.Cells.Item(15, 12).Value = 800
$ws.Cells.Item(15, 13).Value = 750
$ws.Cells.Item(15, 16).Value = 750

$ws.Cells.Item(16, 4).Value = 44764
$ws.Cells.Item(16, 9).Value = "Segunda"
$ws.Cells.Item(16, 10).Value = 150
$ws.Cells.Item(16, 11).Value = 600
$ws.Cells.Item(16, 12).Value = 600
$ws.Cells.Item(16, 13).Value = 600
$ws.Cells.Item(16, 16).Value = 600

$ws.Cells.Item(17, 4).Value = 44791
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 240
$ws.Cells.Item(17, 11).Value = 750
$ws.Cells.Item(17, 12).Value = 800
$ws.Cells.Item(17, 13).Value = 775
$ws.Cells.Item(17, 16).Value = 775

$ws.Cells.Item(18, 4).Value = 44791
$ws.Cells.Item(18, 9).Value = "Segunda"
$ws.Cells.Item(18, 10).Value = 250
$ws.Cells.Item(18, 11).Value = 650
$ws.Cells.Item(18, 12).Value = 650
$ws.Cells.Item(18, 13).Value = 650
$ws.Cells.Item(18, 16).Value = 650

$ws.Cells.Item(19, 4).Value = 45021
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 200
$ws.Cells.Item(19, 11).Value = 1500
$ws.Cells.Item(19, 12).Value = 1500
$ws.Cells.Item(19, 13).Value = 1500
$ws.Cells.Item(19, 16).Value = 1500

$ws.Cells.Item(20, 4).Value = 45030
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 300
$ws.Cells.Item(20, 11).Value = 1500
$ws.Cells.Item(20, 12).Value = 1500
$ws.Cells.Item(20, 13).Value = 1500
$ws.Cells.Item(20, 16).Value = 1500

$ws.Cells.Item(21, 4).Value = 44754
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 200
$ws.Cells.Item(21, 11).Value = 700
$ws.Cells.Item(21, 12).Value = 750
$ws.Cells.Item(21, 13).Value = 725
$ws.Cells.Item(21, 16).Value = 725

$ws.Cells.Item(22, 4).Value = 44533
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 2000
$ws.Cells.Item(22, 12).Value = 2200
$ws.Cells.Item(22, 13).Value = 2100
$ws.Cells.Item(22, 16).Value = 2100

$ws.Cells.Item(23, 4).Value = 44837
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 200
$ws.Cells.Item(23, 11).Value = 700
$ws.Cells.Item(23, 12).Value = 800
$ws.Cells.Item(23, 13).Value = 750
$ws.Cells.Item(23, 16).Value = 750

$ws.Cells.Item(24, 4).Value = 44837
$ws.Cells.Item(24, 9).Value = "Segunda"
$ws.Cells.Item(24, 10).Value = 150
$ws.Cells.Item(24, 11).Value = 600
$ws.Cells.Item(24, 12).Value = 600
$ws.Cells.Item(24, 13).Value = 600
$ws.Cells.Item(24, 16).Value = 600

$ws.Cells.Item(25, 4).Value = 45041
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 200
$ws.Cells.Item(25, 11).Value = 1500
$ws.Cells.Item(25, 12).Value = 1500
$ws.Cells.Item(25, 13).Value = 1500
$ws.Cells.Item(25, 16).Value = 1500

$ws.Cells.Item(26, 4).Value = 45016
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 1500
$ws.Cells.Item(26, 12).Value = 1500
$ws.Cells.Item(26, 13).Value = 1500
$ws.Cells.Item(26, 16).Value = 1500

$ws.Cells.Item(27, 4).Value = 44811
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 300
$ws.Cells.Item(27, 11).Value = 750
$ws.Cells.Item(27, 12).Value = 850
$ws.Cells.Item(27, 13).Value = 800
$ws.Cells.Item(27, 16).Value = 800

$ws.Cells.Item(28, 4).Value = 44831
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 300
$ws.Cells.Item(28, 11).Value = 700
$ws.Cells.Item(28, 12).Value = 800
$ws.Cells.Item(28, 13).Value = 750
$ws.Cells.Item(28, 16).Value = 750

$ws.Cells.Item(29, 4).Value = 44831
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 200
$ws.Cells.Item(29, 11).Value = 600
$ws.Cells.Item(29, 12).Value = 600
$ws.Cells.Item(29, 13).Value = 600
$ws.Cells.Item(29, 16).Value = 600

$ws.Cells.Item(30, 4).Value = 45036
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 250
$ws.Cells.Item(30, 11).Value = 1500
$ws.Cells.Item(30, 12).Value = 1500
$ws.Cells.Item(30, 13).Value = 1500
$ws.Cells.Item(30, 16).Value = 1500

$ws.Cells.Item(31, 4).Value = 44839
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 240
$ws.Cells.Item(31, 11).Value = 700
$ws.Cells.Item(31, 12).Value = 800
$ws.Cells.Item(31, 13).Value = 750
$ws.Cells.Item(31, 16).Value = 750

$ws.Cells.Item(32, 4).Value = 44839
$ws.Cells.Item(32, 9).Value = "Segunda"
$ws.Cells.Item(32, 10).Value = 200
$ws.Cells.Item(32, 11).Value = 600
$ws.Cells.Item(32, 12).Value = 600
$ws.Cells.Item(32, 13).Value = 600
$ws.Cells.Item(32, 16).Value = 600

$ws.Cells.Item(33, 4).Value = 45035
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 150
$ws.Cells.Item(33, 11).Value = 1500
$ws.Cells.Item(33, 12).Value = 1500
$ws.Cells.Item(33, 13).Value = 1500
$ws.Cells.Item(33, 16).Value = 1500

$ws.Cells.Item(34, 4).Value = 45022
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 230
$ws.Cells.Item(34, 11).Value = 1400
$ws.Cells.Item(34, 12).Value = 1500
$ws.Cells.Item(34, 13).Value = 1465
$ws.Cells.Item(34, 16).Value = 1465

$ws.Cells.Item(37, 4).Value = 44859
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 300
$ws.Cells.Item(37, 11).Value = 700
$ws.Cells.Item(37, 12).Value = 800
$ws.Cells.Item(37, 13).Value = 750
$ws.Cells.Item(37, 16).Value = 750

$ws.Cells.Item(38, 4).Value = 44859
$ws.Cells.Item(38, 9).Value = "Segunda"
$ws.Cells.Item(38, 10).Value = 200
$ws.Cells.Item(38, 11).Value = 600
$ws.Cells.Item(38, 12).Value = 600
$ws.Cells.Item(38, 13).Value = 600
$ws.Cells.Item(38, 16).Value = 600

$ws.Cells.Item(39, 4).Value = 44804
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 200
$ws.Cells.Item(39, 11).Value = 750
$ws.Cells.Item(39, 12).Value = 850
$ws.Cells.Item(39, 13).Value = 800
$ws.Cells.Item(39, 16).Value = 800

$ws.Cells.Item(40, 4).Value = 44804
$ws.Cells.Item(40, 9).Value = "Segunda"
$ws.Cells.Item(40, 10).Value = 200
$ws.Cells.Item(40, 11).Value = 650
$ws.Cells.Item(40, 12).Value = 650
$ws.Cells.Item(40, 13).Value = 650
$ws.Cells.Item(40, 16).Value = 650

$ws.Cells.Item(41, 4).Value = 45027
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 200
$ws.Cells.Item(41, 11).Value = 1500
$ws.Cells.Item(41, 12).Value = 1500
$ws.Cells.Item(41, 13).Value = 1500
$ws.Cells.Item(41, 16).Value = 1500

$ws.Cells.Item(42, 4).Value = 44818
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 300
$ws.Cells.Item(42, 11).Value = 800
$ws.Cells.Item(42, 12).Value = 900
$ws.Cells.Item(42, 13).Value = 850
$ws.Cells.Item(42, 16).Value = 850

$ws.Cells.Item(43, 4).Value = 44610
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 100
$ws.Cells.Item(43, 11).Value = 600
$ws.Cells.Item(43, 12).Value = 650
$ws.Cells.Item(43, 13).Value = 625
$ws.Cells.Item(43, 16).Value = 625

$ws.Cells.Item(44, 4).Value = 44608
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 120
$ws.Cells.Item(44, 11).Value = 600
$ws.Cells.Item(44, 12).Value = 650
$ws.Cells.Item(44, 13).Value = 625
$ws.Cells.Item(44, 16).Value = 625

$ws.Cells.Item(45, 4).Value = 44799
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 160
$ws.Cells.Item(45, 11).Value = 750
$ws.Cells.Item(45, 12).Value = 850
$ws.Cells.Item(45, 13).Value = 800
$ws.Cells.Item(45, 16).Value = 800

$ws.Cells.Item(46, 4).Value = 44799
$ws.Cells.Item(46, 9).Value = "Segunda"
$ws.Cells.Item(46, 10).Value = 120
$ws.Cells.Item(46, 11).Value = 650
$ws.Cells.Item(46, 12).Value = 650
$ws.Cells.Item(46, 13).Value = 650
$ws.Cells.Item(46, 16).Value = 650

$ws.Cells.Item(47, 4).Value = 44624
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 120
$ws.Cells.Item(47, 11).Value = 650
$ws.Cells.Item(47, 12).Value = 700
$ws.Cells.Item(47, 13).Value = 675
$ws.Cells.Item(47, 16).Value = 675

$ws.Cells.Item(48, 4).Value = 44883
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 300
$ws.Cells.Item(48, 11).Value = 700
$ws.Cells.Item(48, 12).Value = 800
$ws.Cells.Item(48, 13).Value = 750
$ws.Cells.Item(48, 16).Value = 750

$ws.Cells.Item(49, 4).Value = 44883
$ws.Cells.Item(49, 9).Value = "Segunda"
$ws.Cells.Item(49, 10).Value = 200
$ws.Cells.Item(49, 11).Value = 600
$ws.Cells.Item(49, 12).Value = 600
$ws.Cells.Item(49, 13).Value = 600
$ws.Cells.Item(49, 16).Value = 600
